# Re-sort the calibration data rows (A2:D12) by the "time (s)" column (A)
# in ascending order. The header row (row 1) is left untouched; values are
# written with [double]"..." casts so the exponent-notation literals
# (e.g. "-3.2422421282e-06") parse correctly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"57034.363191"
$ws.Cells.Item(2, 2).Value = [double]"-3.2422421282e-06"
$ws.Cells.Item(2, 3).Value = [double]"-8.707349760600001e-06"
$ws.Cells.Item(2, 4).Value = [double]"-2.4915333007e-05"
$ws.Cells.Item(3, 1).Value = [double]"57046.031192"
$ws.Cells.Item(3, 2).Value = [double]"-6.6157834425e-05"
$ws.Cells.Item(3, 3).Value = [double]"-9.6710380647e-05"
$ws.Cells.Item(3, 4).Value = [double]"-0.00010499781355"
$ws.Cells.Item(4, 1).Value = [double]"57059.299192"
$ws.Cells.Item(4, 2).Value = [double]"-7.372849520099999e-05"
$ws.Cells.Item(4, 3).Value = [double]"-0.00022851370616"
$ws.Cells.Item(4, 4).Value = [double]"-0.00015041768697"
$ws.Cells.Item(5, 1).Value = [double]"57070.563193"
$ws.Cells.Item(5, 2).Value = [double]"-8.153237473799999e-05"
$ws.Cells.Item(5, 3).Value = [double]"-0.00036044991622"
$ws.Cells.Item(5, 4).Value = [double]"-0.0001933029678"
$ws.Cells.Item(6, 1).Value = [double]"57081.963193"
$ws.Cells.Item(6, 2).Value = [double]"-0.0001228416"
$ws.Cells.Item(6, 3).Value = [double]"-0.0004886953"
$ws.Cells.Item(6, 4).Value = [double]"-0.0002440999"
$ws.Cells.Item(7, 1).Value = [double]"57091.831194"
$ws.Cells.Item(7, 2).Value = [double]"-0.0001927213"
$ws.Cells.Item(7, 3).Value = [double]"-0.0006102247000000001"
$ws.Cells.Item(7, 4).Value = [double]"-0.0003006622"
$ws.Cells.Item(8, 1).Value = [double]"57102.299194"
$ws.Cells.Item(8, 2).Value = [double]"-0.000146344"
$ws.Cells.Item(8, 3).Value = [double]"-0.0004822301"
$ws.Cells.Item(8, 4).Value = [double]"-0.0002532177"
$ws.Cells.Item(9, 1).Value = [double]"57112.163195"
$ws.Cells.Item(9, 2).Value = [double]"-9.6978917629e-05"
$ws.Cells.Item(9, 3).Value = [double]"-0.00035120831059"
$ws.Cells.Item(9, 4).Value = [double]"-0.00020319603157"
$ws.Cells.Item(10, 1).Value = [double]"57123.231195"
$ws.Cells.Item(10, 2).Value = [double]"-4.4424733057e-05"
$ws.Cells.Item(10, 3).Value = [double]"-0.00022768781771"
$ws.Cells.Item(10, 4).Value = [double]"-0.00014148830698"
$ws.Cells.Item(11, 1).Value = [double]"57136.031196"
$ws.Cells.Item(11, 2).Value = [double]"-2.6573538511e-05"
$ws.Cells.Item(11, 3).Value = [double]"-7.3390219873e-05"
$ws.Cells.Item(11, 4).Value = [double]"-0.00010844868675"
$ws.Cells.Item(12, 1).Value = [double]"57147.963197"
$ws.Cells.Item(12, 2).Value = [double]"-9.245939921800001e-06"
$ws.Cells.Item(12, 3).Value = [double]"-1.3298685729e-05"
$ws.Cells.Item(12, 4).Value = [double]"-2.5306092296e-05"
